$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: update existing K10 value
$ws.Range("K10").Value = 1.388747888886706

# Row 11: update existing J11, add new K11
$ws.Range("J11").Value = 0.444773652920949
$ws.Range("K11").Value = 0.2348700177716323

# Row 12: update existing I12, add new J12
$ws.Range("I12").Value = 0.4487415504340581
$ws.Range("J12").Value = 0.2388379152847414

# Row 13: update existing H13, add new I13
$ws.Range("H13").Value = 0.5843816406042994
$ws.Range("I13").Value = 0.3744780054549828

# Row 14: update existing G14, add new H14
$ws.Range("G14").Value = 0.3435754587486348
$ws.Range("H14").Value = 0.1336718235993181

# Row 15: update existing F15, add new G15
$ws.Range("F15").Value = 0.2982442434965384
$ws.Range("G15").Value = 0.08834060834722172

# Row 16: update existing E16, add new F16
$ws.Range("E16").Value = 0.2313828215604846
$ws.Range("F16").Value = 0.02147918641116785

# Row 17: update existing D17, add new E17
$ws.Range("D17").Value = 0.201796619203768
$ws.Range("E17").Value = -0.00810701594554874

# Row 18: update existing C18, add new D18
$ws.Range("C18").Value = 0.1836459624741271
$ws.Range("D18").Value = -0.02625767267518964

# Row 19: update existing B19, add new C19
$ws.Range("B19").Value = 0.1656141382254278
$ws.Range("C19").Value = -0.04428949692388896

# Row 20: add new B20
$ws.Range("B20").Value = -0.09587373626955231
